$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (column B) only had a header ("Price") with no data
# underneath it. Delete the whole column outright: this shifts the
# FirstName/LastName columns left into B/C and drops the now-unused
# "Price" entry from the shared strings table.
$ws.Columns("B").Delete()

# Reflect the post-edit selection state: column B (now FirstName) selected
# with B1 as the active cell.
$ws.Range("B1:B1048576").Select()
